$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.820.55'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").Value = '2.587.20'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.08%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").Value = '2.594.72'
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.101'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("E13").Value = '  +2.98%  '
$ws.Range("D14").Value = '3.039.72'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("D15").Value = '58.722.79'
$ws.Range("E15").Value = '  +2.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.31%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000133'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.574.37'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '339.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.76%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.168'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '0.0₃0720'
$ws.Range("E30").Value = '  -3.89%  '
$ws.Range("E31").Value = '  -4.73%  '
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.97'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.54'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.47'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.829'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.816'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.51'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '275.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.67%  '
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.590'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("D49").Value = '1.982.67'
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.57%  '
